# Update column C values (rows 2-12) on Sheet1 per the commit's bug fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 0.7621
    3  = 0.9043
    4  = 0.6504
    5  = 0.9034
    6  = 0.8723
    7  = 0.5569
    8  = 0.6069
    9  = 0.7092000000000001
    10 = 0.7276
    11 = 0.7517
    12 = 0.4752
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
